# Insert a new weekly record before the existing row 29 ("Feria Lagunitas de
# Puerto Montt" / Sandia sheet). This pushes the previous rows 29-141 down to
# 30-142 (dimension grows from A1:R141 to A1:R142) and fills the newly
# inserted row 29 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29, shifting rows 29:141 down to 30:142.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new record.
$ws.Range("A29").Value = 4
$ws.Range("B29").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C29").Value = "Los Lagos"
$ws.Range("D29").Value = 44525
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 100112028
$ws.Range("G29").Value = "Sandia"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 600
$ws.Range("K29").Value = 1200
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = 1200
$ws.Range("N29").Value = "$/kilo (volumen en unidades)"
$ws.Range("O29").Value = "Perú"
$ws.Range("P29").Value = 1200
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
